$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 234356.48
$ws.Range("I113").Value = 385738.3
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 385738.3
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = -382484.3
$ws.Range("N113").Value = -9808
# Row 116
$ws.Range("H116").Value = 3254307.8
$ws.Range("I116").Value = 17859192
$ws.Range("J116").Value = 8777.777
$ws.Range("K116").Value = 17859192
$ws.Range("L116").Value = 8777.777
$ws.Range("M116").Value = -17855750
$ws.Range("N116").Value = -15661.777
# Row 138
$ws.Range("H138").Value = 2578.411
$ws.Range("I138").Value = 877.7692
$ws.Range("J138").Value = 4529.147
$ws.Range("K138").Value = 2633.3076
$ws.Range("L138").Value = 13587.441
$ws.Range("M138").Value = 2506.6924
$ws.Range("N138").Value = -23867.441

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 5560.7144
$ws.Range("I2").Value = 8481.25
$ws.Range("J2").Value = 1666.6666
$ws.Range("K2").Value = 8481.25
$ws.Range("L2").Value = 1666.6666
$ws.Range("M2").Value = -8368.25
$ws.Range("N2").Value = -1892.6666
# Row 45
$ws.Range("H45").Value = 2050.5
$ws.Range("I45").Value = 1601
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1601
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -1224
$ws.Range("N45").Value = -3254
# Row 61
$ws.Range("H61").Value = 1538.6818
$ws.Range("I61").Value = 1249.4
$ws.Range("J61").Value = 2158.5715
$ws.Range("K61").Value = 1249.4
$ws.Range("L61").Value = 2158.5715
$ws.Range("M61").Value = -1037.4
$ws.Range("N61").Value = -2582.5715
# Row 63
$ws.Range("H63").Value = 2356.254
$ws.Range("I63").Value = 2331.9656
$ws.Range("J63").Value = 2638
$ws.Range("K63").Value = 2331.9656
$ws.Range("L63").Value = 2638
$ws.Range("M63").Value = -1645.9656
$ws.Range("N63").Value = -4010
# Row 66
$ws.Range("H66").Value = 2356.254
$ws.Range("I66").Value = 2331.9656
$ws.Range("J66").Value = 2638
$ws.Range("K66").Value = 11659.828
$ws.Range("L66").Value = 13190
$ws.Range("M66").Value = -8227.828
$ws.Range("N66").Value = -20054
# Row 116
$ws.Range("H116").Value = 5560.7144
$ws.Range("I116").Value = 8481.25
$ws.Range("J116").Value = 1666.6666
$ws.Range("K116").Value = 8481.25
$ws.Range("L116").Value = 1666.6666
$ws.Range("M116").Value = -6187.25
$ws.Range("N116").Value = -6254.6666
# Row 136
$ws.Range("H136").Value = 1538.6818
$ws.Range("I136").Value = 1249.4
$ws.Range("J136").Value = 2158.5715
$ws.Range("K136").Value = 3748.2
$ws.Range("L136").Value = 6475.7145
$ws.Range("M136").Value = -1198.2
$ws.Range("N136").Value = -11575.7145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 5560.7144
$ws.Range("I3").Value = 8481.25
$ws.Range("J3").Value = 1666.6666
$ws.Range("K3").Value = 8481.25
$ws.Range("L3").Value = 1666.6666
$ws.Range("M3").Value = -8367.25
$ws.Range("N3").Value = -1894.6666
# Row 20
$ws.Range("H20").Value = 16134641
$ws.Range("I20").Value = 22733816
$ws.Range("J20").Value = 3323.889
$ws.Range("K20").Value = 22733816
$ws.Range("L20").Value = 3323.889
$ws.Range("M20").Value = -22733569
$ws.Range("N20").Value = -3817.889
# Row 55
$ws.Range("H55").Value = 69175
$ws.Range("J55").Value = 69175
$ws.Range("L55").Value = 69175
$ws.Range("N55").Value = -69721

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 829.7368
$ws.Range("I58").Value = 622.1212
$ws.Range("J58").Value = 2200
$ws.Range("K58").Value = 622.1212
$ws.Range("L58").Value = 2200
$ws.Range("M58").Value = -419.1212
$ws.Range("N58").Value = -2606
# Row 132
$ws.Range("H132").Value = 7494.3335
$ws.Range("I132").Value = 8305.1875
$ws.Range("J132").Value = 4899.6
$ws.Range("K132").Value = 24915.5625
$ws.Range("L132").Value = 14698.8
$ws.Range("M132").Value = -22385.5625
$ws.Range("N132").Value = -19758.8
# Row 136
$ws.Range("H136").Value = 829.7368
$ws.Range("I136").Value = 622.1212
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 1866.3636
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = 683.6363999999999
$ws.Range("N136").Value = -11700

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 60208.824
$ws.Range("I4").Value = 182
$ws.Range("J4").Value = 340334
$ws.Range("K4").Value = 546
$ws.Range("L4").Value = 1021002
$ws.Range("M4").Value = -434
$ws.Range("N4").Value = -1021226
# Row 40
$ws.Range("H40").Value = 347.8
$ws.Range("I40").Value = 133.42857
$ws.Range("J40").Value = 848
$ws.Range("K40").Value = 533.71428
$ws.Range("L40").Value = 3392
$ws.Range("M40").Value = -464.71428
$ws.Range("N40").Value = -3530
# Row 51
$ws.Range("H51").Value = 2441.4707
$ws.Range("J51").Value = 2885
$ws.Range("L51").Value = 8655
$ws.Range("N51").Value = -9575
# Row 103
$ws.Range("H103").Value = 2128
$ws.Range("I103").Value = 342.5
$ws.Range("J103").Value = 2265.3462
$ws.Range("K103").Value = 1027.5
$ws.Range("L103").Value = 6796.0386
$ws.Range("M103").Value = -148.5
$ws.Range("N103").Value = -8554.0386
# Row 121
$ws.Range("H121").Value = 26319298
$ws.Range("J121").Value = 27781454
$ws.Range("L121").Value = 83344362
$ws.Range("N121").Value = -83346982

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null
# Row 102
$ws.Range("H102").Value = 2150.4167
$ws.Range("I102").Value = 1468
$ws.Range("K102").Value = 1468
$ws.Range("M102").Value = 154
# Row 122
$ws.Range("H122").Value = 2962.5
$ws.Range("I122").Value = 2600
$ws.Range("J122").Value = 3566.6667
$ws.Range("K122").Value = 7800
$ws.Range("L122").Value = 10700.0001
$ws.Range("M122").Value = -5350
$ws.Range("N122").Value = -15600.0001
# Row 132
$ws.Range("H132").Value = 5751.769
$ws.Range("I132").Value = 6487.4
$ws.Range("J132").Value = 3299.6667
$ws.Range("K132").Value = 19462.2
$ws.Range("L132").Value = 9899.000100000001
$ws.Range("M132").Value = -16932.2
$ws.Range("N132").Value = -14959.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2095.4614
$ws.Range("I40").Value = 2014.1
$ws.Range("J40").Value = 2366.6667
$ws.Range("K40").Value = 2014.1
$ws.Range("L40").Value = 2366.6667
$ws.Range("M40").Value = -1878.1
$ws.Range("N40").Value = -2638.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3832.1667
$ws.Range("I62").Value = 3990
$ws.Range("J62").Value = 3800.6
$ws.Range("K62").Value = 3990
$ws.Range("L62").Value = 3800.6
$ws.Range("M62").Value = -3366
$ws.Range("N62").Value = -5048.6
# Row 65
$ws.Range("H65").Value = 3832.1667
$ws.Range("I65").Value = 3990
$ws.Range("J65").Value = 3800.6
$ws.Range("K65").Value = 19950
$ws.Range("L65").Value = 19003
$ws.Range("M65").Value = -16830
$ws.Range("N65").Value = -25243
# Row 96
$ws.Range("H96").Value = 2395.5
$ws.Range("I96").Value = 2286.5
$ws.Range("J96").Value = 2450
$ws.Range("K96").Value = 2286.5
$ws.Range("L96").Value = 2450
$ws.Range("M96").Value = -913.5
$ws.Range("N96").Value = -5196
# Row 100
$ws.Range("H100").Value = 11364805
$ws.Range("I100").Value = 15152915
$ws.Range("J100").Value = 475
$ws.Range("K100").Value = 30305830
$ws.Range("L100").Value = 950
$ws.Range("M100").Value = -30305289
$ws.Range("N100").Value = -2032
# Row 107
$ws.Range("H107").Value = 6667.4375
$ws.Range("I107").Value = 409.125
$ws.Range("J107").Value = 12925.75
$ws.Range("K107").Value = 1227.375
$ws.Range("L107").Value = 38777.25
$ws.Range("M107").Value = 692.625
$ws.Range("N107").Value = -42617.25
